# Fruta / hortaliza, semanal
# A new daily price record was added to the top of the data table
# (first data row, row 5), pushing all the existing records down by
# one row. The previously-last record (old row 101) now lives in the
# newly created row 102.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 5:101 down to 6:102, leaving row 5 empty.
$ws.Rows("5:5").Insert()

# Populate the newly inserted row 5 with the new record.
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "Vega Modelo de Temuco"
$ws.Range("C5").Value = "La Araucanía"
$ws.Range("D5").Value = 45092
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100108
$ws.Range("H5").Value = "Tropicales y subtropicales"
$ws.Range("I5").Value = 100108003
$ws.Range("J5").Value = "Maracuyá"
$ws.Range("K5").Value = "Sin especificar"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 50000
$ws.Range("O5").Value = 50000
$ws.Range("P5").Value = 50000
$ws.Range("Q5").Value = "$/caja 18 kilos"
$ws.Range("R5").Value = "Región de Arica y Parinacota"
$ws.Range("S5").Value = 2778
$ws.Range("T5").Value = 18
